$d = $word.ActiveDocument

# Remove the stray "_GoBack" bookmark (w:bookmarkStart/w:bookmarkEnd) left in
# the first paragraph.
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    if ($bm -ne $null) {
        $bm.Delete()
    }
} catch {
    # no-op: bookmark not present
}

# Remove the empty paragraph, the bold "Как комплементарные цвета влияют на
# наше восприятие?" question paragraph, the empty paragraph after it, and the
# "Они создают баланс..." answer paragraph that follows the first Q&A pair.
$start = $d.Paragraphs.Item(4)
$end = $d.Paragraphs.Item(7)
$r = $d.Range($start.Range.Start, $end.Range.End)
$r.Delete()
